$d = $word.ActiveDocument

# --- Locate the end of the sentence "...getting ate in his absence." ---
# (the paragraph that currently holds the _GoBack bookmark at its end).
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "getting ate in his absence.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor sentence to split the paragraph."
}
$splitPoint = $findRange.End

# --- Split the paragraph in two, right after that sentence. ---
# A collapsed range at the split point is used so the new paragraph
# inherits the same style/numbering (ListParagraph) as the original one,
# and the existing "_GoBack" bookmark (sitting right at the split point)
# stays attached to the end of the first half.
$splitRng = $d.Range($splitPoint, $splitPoint)
$splitRng.InsertParagraphAfter()

# --- The new (second) paragraph is the one that now follows. ---
$newPara = $d.Paragraphs(2).Next()

# --- Fill it with the full constraint text. ---
$constraintsText = "The constraints of the problem are that the man cannot leave any of the two behind without something going wrong. He is only allowed to bring one thing at a time across."
$newPara.Range.InsertBefore($constraintsText)

# --- Move the "_GoBack" bookmark so it sits inside the new sentence, ---
# --- right after "The constraints of the problem are t".            ---
$prefix = "The constraints of the problem are t"
$newParaStart = $newPara.Range.Start
$bookmarkPos = $newParaStart + $prefix.Length

$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
